$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("E10").Value = "-"

# Row 11
$ws.Range("E11").Value = "['MEC-2A-Des. Maq. Cad_T2', 'MEC-2A-Des. Maq. Cad_T2']"

# Row 12
$ws.Range("E12").Value = "['MEC-2A-Des. Maq. Cad_T2', -]"

# Row 16
$ws.Range("E16").Value = "-"

# Row 18
$ws.Range("D18").Value = "-"

# Row 19
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "MEC-2NB-Elemaq."

# Row 20
$ws.Range("C20").Value = "MEC-2NA-Des. Maq. Cad"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "MEC-2NB-Elemaq."
$ws.Range("F20").Value = "MEC-2NA-Des. Maq. Cad"

# Row 21
$ws.Range("C21").Value = "MEC-2NA-Des. Maq. Cad"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "MEC-2NA-Elemaq."
